$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
Write-Output $win
$win.TabRatio = 500
Write-Output "set"
